$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Helper: insert a brand-new, clean empty paragraph right after $afterPara
# (no stray placeholder run left behind) and return the new Paragraph.
# ----------------------------------------------------------------------
function New-CleanParagraphAfter($afterPara) {
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $idx = $afterPara.Index + 1
    $newPara = $d.Paragraphs.Item($idx)
    $nr = $newPara.Range
    $nr.Text = "x"
    $freshR = $d.Paragraphs.Item($idx).Range
    $freshR.MoveEnd(1, -1) | Out-Null
    $freshR.Delete()
    return $d.Paragraphs.Item($idx)
}

# ----------------------------------------------------------------------
# Helper: insert a brand-new paragraph right after $afterPara whose whole
# content is a hyperlink with the given address/display text. Returns the
# new Paragraph.
# ----------------------------------------------------------------------
function New-HyperlinkParagraphAfter($afterPara, [string]$url) {
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $idx = $afterPara.Index + 1
    $newPara = $d.Paragraphs.Item($idx)
    $nr = $newPara.Range
    $nr.Text = $url
    $freshR = $d.Paragraphs.Item($idx).Range
    $freshR.MoveEnd(1, -1) | Out-Null
    $d.Hyperlinks.Add($freshR, $url, "", "", $url) | Out-Null
    return $d.Paragraphs.Item($idx)
}

# Paragraph holding the "http://www.stat.columbia.edu/~gelman/book/data/" link.
$gelmanPara = $d.Paragraphs.Item(11)

# Add the two new hyperlink paragraphs right after it.
$ucla = New-HyperlinkParagraphAfter $gelmanPara "http://www.ats.ucla.edu/stat/r/dae/logit.htm"
$logistic = New-HyperlinkParagraphAfter $ucla "http://logisticregressionanalysis.com/303-what-a-logistic-regression-data-set-looks-like-an-example/"

# The paragraph that follows the new hyperlinks is the _GoBack bookmark paragraph.
$bookmarkPara = $d.Paragraphs.Item($logistic.Index + 1)

# Add three new blank paragraphs right after the bookmark paragraph.
$p = $bookmarkPara
for ($i = 0; $i -lt 3; $i++) {
    $p = New-CleanParagraphAfter $p
}
